$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.148.07"
$ws.Range("E2").Value = "  -4.39%  "
$ws.Range("D3").Value = "1.652.99"
$ws.Range("E3").Value = "  -3.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.10"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5108"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  -3.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06424"
$ws.Range("E9").Value = "  -3.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.99"
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "1.655.08"
$ws.Range("E12").Value = "  -3.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.280"
$ws.Range("E13").Value = "  -4.95%  "
$ws.Range("D14").Value = "1.880.45"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5519"
$ws.Range("E15").Value = "  -5.45%  "
$ws.Range("D16").Value = "0.0₅8031"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.03"
$ws.Range("E17").Value = "  -5.85%  "
$ws.Range("D18").Value = "26.141.68"
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.70"
$ws.Range("E20").Value = "  -5.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.412"
$ws.Range("E21").Value = "  -4.71%  "
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.035"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.51"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.737"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1181"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.979"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.86"
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05112"
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.245"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.343"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.223"
$ws.Range("E33").Value = "  -6.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.566"
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.361"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9256"
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").Value = "1.171.15"
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5687"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01589"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8315"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.656"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.33"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "1.790.00"
$ws.Range("E46").Value = "  -3.47%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4548"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.66"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.876"
$ws.Range("E51").Value = "  -2.83%  "
